## "added prob 19 in common part"
##
## The workbook's Sheet1 lists building-block pairs (code / description /
## math snippet) for a set of problems. This edit:
##   1. Shortens the wording of the existing B45 entry (drops the redundant
##      "<sum> 성질을 이용해서" lead-in, keeping just the second clause).
##   2. Appends two new rows (77 and 78) for "problem 19" in the common
##      part: a derivative-inequality step and a discriminant-inequality
##      step, each with a code (column A), a description (column B) and a
##      math expression (column C).
##   3. Moves the sheet's scroll position / active selection down near the
##      newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Shorten the existing B45 description.
$ws.Range("B45").Value = '$\displaystyle\sum$가 포함된 두 식을 연립합니다.'

# 2) Add the two new rows for problem 19 (codes + descriptions first,
#    then the math snippets, matching how the rows were authored).
$ws.Range("A77").Value = "d0027"
$ws.Range("B77").Value = '도함수에 대한 부등식을 세웁니다.'

$ws.Range("A78").Value = "d0028"
$ws.Range("B78").Value = '이차대부등식이 항상 성립하도록 판별식에 대한 부등식을 세웁니다.'

$ws.Range("C77").Value = '$f^{\prime}(x) \geq 0$;'
$ws.Range("C78").Value = '$\dfrac{D}{4} \leq 0$;'

# 3) Scroll the view toward the new rows and move the active selection.
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A79").Select()
